# Weekly update: insert a new record at the top of the Zanahoria (carrot)
# price series (row 149), pushing all subsequent rows down by one.
# The last existing row (previously 232) becomes row 233.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 149; this shifts rows 149..232 down to
# 150..233 and extends the worksheet dimension to A1:R233 automatically.
$ws.Rows.Item(149).Insert()

# Populate the newly inserted row 149 with the new weekly record.
$ws.Cells.Item(149, 1).Value  = 10
$ws.Cells.Item(149, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(149, 3).Value  = "La Araucanía"
$ws.Cells.Item(149, 4).Value  = 44529
$ws.Cells.Item(149, 5).Value  = 9
$ws.Cells.Item(149, 6).Value  = 100114013
$ws.Cells.Item(149, 7).Value  = "Zanahoria"
$ws.Cells.Item(149, 8).Value  = "Sin especificar"
$ws.Cells.Item(149, 9).Value  = "Primera"
$ws.Cells.Item(149, 10).Value = 185
$ws.Cells.Item(149, 11).Value = 6000
$ws.Cells.Item(149, 12).Value = 7000
$ws.Cells.Item(149, 13).Value = 6676
$ws.Cells.Item(149, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(149, 15).Value = "Región del Maule"
$ws.Cells.Item(149, 16).Value = 334
$ws.Cells.Item(149, 17).Value = 20
$ws.Cells.Item(149, 18).Value = "Hortaliza"
